{"js": "// Bump the \"Version 1.\" label to \"Version 2.\" (wireframes revert-of-a-revert).\n// Net visible-text change per the diff: \"Version 1.\" -> \"Version 2.\"\nconst body = context.document.body;\n\n// Scope the search tightly to \" 1.\" so we only touch the version number/\n// trailing period, not any other lone \"1\" that might appear in the doc.\nconst results = body.search(\" 1.\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\" 2.\", Word.InsertLocation.replace);\n} else {\n  // Fallback: the exact \" 1.\" run wasn't found (e.g. formatting already\n  // differs) -- try the bare digit next to \"Version\" instead.\n  const fallback = body.search(\"Version 1\", { matchCase: true });\n  fallback.load(\"text\");\n  await context.sync();\n  if (fallback.items.length > 0) {\n    fallback.items[0].insertText(\"Version 2\", Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Bump the \"Version 1.\" label to \"Version 2.\" (wireframes revert-of-a-revert).\n# Net visible-text change per the diff: \"Version 1.\" -> \"Version 2.\"\n$d = $word.ActiveDocument\n\n# Scope the search tightly to \" 1.\" so we only touch the version number and\n# its trailing period, not any other lone \"1\" that might appear in the doc.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \" 1.\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \" 2.\"\n$found = $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\nif (-not $found) {\n    # Fallback: the exact \" 1.\" run wasn't found (e.g. formatting already\n    # differs) -- try the looser \"Version 1\" phrase instead.\n    $find2 = $d.Content.Find\n    $find2.ClearFormatting()\n    $find2.Text = \"Version 1\"\n    $find2.Replacement.ClearFormatting()\n    $find2.Replacement.Text = \"Version 2\"\n    $find2.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n}\n"}
